$p = $ppt.ActivePresentation
Write-Host "Before: $($p.HasNotesMaster)"
try {
  $p.HasNotesMaster = $false
  Write-Host "set false OK: $($p.HasNotesMaster)"
} catch { Write-Host "ERR false: $_" }
try {
  $p.HasNotesMaster = $true
  Write-Host "set true OK: $($p.HasNotesMaster)"
} catch { Write-Host "ERR true: $_" }
